$wb = $excel.ActiveWorkbook

# --- Update the publication Date on the "Metadata" sheet ---
$wsMeta = $wb.Worksheets.Item(1)
$wsMeta.Range("B8").Value = "2023-02-21T11:59:56+00:00"

# --- Append four new concept rows to the "Concepts" sheet ---
$ws = $wb.Worksheets.Item(2)

# Copy the formatting (borders / wrap / style) of the last existing data
# row (row 8) down into the four new rows (9-12).
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Copy the "Level" value ("1", stored as text) from A2 into each new row
# in column A, so it keeps being a text value instead of turning into a
# number.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4163)
$ws.Range("A10").PasteSpecial(-4163)
$ws.Range("A11").PasteSpecial(-4163)
$ws.Range("A12").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Row 9: BW
$ws.Range("B9").Value = "BW"
$ws.Range("C9").Value = "BW File"

# Row 10: BED
$ws.Range("B10").Value = "BED"
$ws.Range("C10").Value = "BED File"

# Row 11: PNG
$ws.Range("B11").Value = "PNG"
$ws.Range("C11").Value = "PNG File"

# Row 12: CSV
$ws.Range("B12").Value = "CSV"
$ws.Range("C12").Value = "CSV File"
